# Add a new sample record (row 3) to the password sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "/B0talDhaniya4"
$ws.Range("B3").Value = "Testing"
